$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 18.63888133333333
$ws.Range("H2").Value2 = 55.91664400000001
$ws.Range("I2").Value2 = 0.5880476104010496
$ws.Range("J2").Value2 = 0.5880476104010497
$ws.Range("M2").Value2 = 6.060959
$ws.Range("N2").Value2 = 18.182877
$ws.Range("O2").Value2 = 0.5093955725568765
$ws.Range("P2").Value2 = 0.5093955725568765
$ws.Range("Q2").Value2 = 112.9694955671987
$ws.Range("R2").Value2 = 1016.725460104788
$ws.Range("S2").Value2 = 0.2995488491909457
$ws.Range("T2").Value2 = 0.2995488491909458

# Row 3
$ws.Range("G3").Value2 = 18.63888133333333
$ws.Range("H3").Value2 = 55.91664400000001
$ws.Range("I3").Value2 = 0.5880476104010496
$ws.Range("J3").Value2 = 0.5880476104010497
$ws.Range("O3").Value2 = 0.4585214502287212
$ws.Range("P3").Value2 = 0.4585214502287212
$ws.Range("Q3").Value2 = 101.6870576221889
$ws.Range("R3").Value2 = 915.1835185997002
$ws.Range("S3").Value2 = 0.2696324431246233
$ws.Range("T3").Value2 = 0.2696324431246233

# Row 4
$ws.Range("G4").Value2 = 18.63888133333333
$ws.Range("H4").Value2 = 55.91664400000001
$ws.Range("I4").Value2 = 0.5880476104010496
$ws.Range("J4").Value2 = 0.5880476104010497
$ws.Range("O4").Value2 = 0.03208297721440233
$ws.Range("P4").Value2 = 0.03208297721440233
$ws.Range("Q4").Value2 = 7.115094726898668
$ws.Range("R4").Value2 = 64.03585254208801
$ws.Range("S4").Value2 = 0.01886631808548061
$ws.Range("T4").Value2 = 0.01886631808548062

# Row 5
$ws.Range("I5").Value2 = 0.3007497405000055
$ws.Range("J5").Value2 = 0.3007497405000055
$ws.Range("M5").Value2 = 6.060959
$ws.Range("N5").Value2 = 18.182877
$ws.Range("O5").Value2 = 0.5093955725568765
$ws.Range("P5").Value2 = 0.5093955725568765
$ws.Range("Q5").Value2 = 57.77686342961265
$ws.Range("R5").Value2 = 519.9917708665139
$ws.Range("S5").Value2 = 0.1532005862583323
$ws.Range("T5").Value2 = 0.1532005862583323

# Row 6
$ws.Range("I6").Value2 = 0.3007497405000055
$ws.Range("J6").Value2 = 0.3007497405000055
$ws.Range("O6").Value2 = 0.4585214502287212
$ws.Range("P6").Value2 = 0.4585214502287212
$ws.Range("S6").Value2 = 0.1379002071699741
$ws.Range("T6").Value2 = 0.1379002071699741

# Row 7
$ws.Range("I7").Value2 = 0.3007497405000055
$ws.Range("J7").Value2 = 0.3007497405000055
$ws.Range("O7").Value2 = 0.03208297721440233
$ws.Range("P7").Value2 = 0.03208297721440233
$ws.Range("S7").Value2 = 0.009648947071699089
$ws.Range("T7").Value2 = 0.009648947071699091

# Row 8
$ws.Range("G8").Value2 = 3.524702666666666
$ws.Range("I8").Value2 = 0.1112026490989449
$ws.Range("J8").Value2 = 0.1112026490989449
$ws.Range("M8").Value2 = 6.060959
$ws.Range("N8").Value2 = 18.182877
$ws.Range("O8").Value2 = 0.5093955725568765
$ws.Range("P8").Value2 = 0.5093955725568765
$ws.Range("Q8").Value2 = 21.36307834985733
$ws.Range("R8").Value2 = 192.2677051487159
$ws.Range("S8").Value2 = 0.05664613710759845
$ws.Range("T8").Value2 = 0.05664613710759846

# Row 9
$ws.Range("G9").Value2 = 3.524702666666666
$ws.Range("I9").Value2 = 0.1112026490989449
$ws.Range("J9").Value2 = 0.1112026490989449
$ws.Range("O9").Value2 = 0.4585214502287212
$ws.Range("P9").Value2 = 0.4585214502287212
$ws.Range("S9").Value2 = 0.05098879993412379
$ws.Range("T9").Value2 = 0.0509887999341238

# Row 10
$ws.Range("G10").Value2 = 3.524702666666666
$ws.Range("I10").Value2 = 0.1112026490989449
$ws.Range("J10").Value2 = 0.1112026490989449
$ws.Range("O10").Value2 = 0.03208297721440233
$ws.Range("P10").Value2 = 0.03208297721440233
$ws.Range("Q10").Value2 = 1.345498847757333
$ws.Range("S10").Value2 = 0.003567712057222626
$ws.Range("T10").Value2 = 0.003567712057222627

$wb.Save()